$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-48
# D column values are forced to Text format to avoid Excel auto-converting
# numeric-looking strings (e.g. "1.00") into actual numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.117.29"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.16"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.71"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.35"
$ws.Range("E7").Value = "  +6.54%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0780"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.82"
$ws.Range("E12").Value = "  +5.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.351.23"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.811"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.53"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.052.58"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.150.70"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.74"
$ws.Range("E18").Value = "  +17.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.63"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.42"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  +10.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.05"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  +9.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.69"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0611"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0898"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.22"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("E38").Value = "  +6.72%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +14.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.10"
$ws.Range("E41").Value = "  +28.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.47"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.12"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.77"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.278.14"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  -1.90%  "
# Row 49/50: the two coins (RocketPoolETH and FraxShare) swapped positions,
# each picking up updated Price/Volume values. Row 51 (FTXToken) price/volume updated too.
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.73"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.240.60"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.40"
$ws.Range("E51").Value = "  -19.43%  "
